$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.416.98"
$ws.Range("E2").Value = "'  -2.18%  "

$ws.Range("D3").Value = "'2.889.19"
$ws.Range("E3").Value = "'  -2.08%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.04%  "

$ws.Range("D5").Value = "'566.17"
$ws.Range("E5").Value = "'  -4.63%  "

$ws.Range("D6").Value = "'142.94"
$ws.Range("E6").Value = "'  -3.38%  "

$ws.Range("E7").Value = "'  +0.12%  "

$ws.Range("E8").Value = "'  -0.53%  "

$ws.Range("D9").Value = "'2.889.77"

$ws.Range("E10").Value = "'  -6.80%  "

$ws.Range("E11").Value = "'  -3.74%  "

$ws.Range("D12").Value = "'0.433"
$ws.Range("E12").Value = "'  -2.31%  "

$ws.Range("E13").Value = "'  -1.94%  "

$ws.Range("D14").Value = "'31.85"
$ws.Range("E14").Value = "'  -3.22%  "

$ws.Range("E15").Value = "'  -0.66%  "

$ws.Range("D16").Value = "'3.367.97"
$ws.Range("E16").Value = "'  -2.08%  "

$ws.Range("D17").Value = "'61.421.77"
$ws.Range("E17").Value = "'  -2.11%  "

$ws.Range("D18").Value = "'6.55"
$ws.Range("E18").Value = "'  -2.57%  "

$ws.Range("D19").Value = "'2.885.89"
$ws.Range("E19").Value = "'  -2.41%  "

$ws.Range("D20").Value = "'431.24"
$ws.Range("E20").Value = "'  -2.61%  "

$ws.Range("D21").Value = "'13.07"
$ws.Range("E21").Value = "'  -3.05%  "

$ws.Range("D22").Value = "'0.653"
$ws.Range("E22").Value = "'  -2.25%  "

$ws.Range("D23").Value = "'6.81"
$ws.Range("E23").Value = "'  -3.20%  "

$ws.Range("D24").Value = "'79.04"
$ws.Range("E24").Value = "'  -2.93%  "

$ws.Range("D25").Value = "'11.81"
$ws.Range("E25").Value = "'  +0.83%  "

$ws.Range("B26").Value = "'Dai"
$ws.Range("C26").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "'  +0.02%  "

$ws.Range("B27").Value = "'RenderToken"
$ws.Range("C27").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'9.97"
$ws.Range("E27").Value = "'  -10.65%  "

$ws.Range("D28").Value = "'2.00"
$ws.Range("E28").Value = "'  -7.02%  "

$ws.Range("E29").Value = "'  -0.50%  "

$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "'  -3.81%  "

$ws.Range("E31").Value = "'  -4.42%  "

$ws.Range("D32").Value = "'2.05"
$ws.Range("E32").Value = "'  -8.62%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "'  +0.08%  "

$ws.Range("E34").Value = "'  -2.88%  "

$ws.Range("D35").Value = "'25.51"
$ws.Range("E35").Value = "'  -3.76%  "

$ws.Range("D36").Value = "'0.955"
$ws.Range("E36").Value = "'  -3.67%  "

$ws.Range("D37").Value = "'5.39"
$ws.Range("E37").Value = "'  -4.03%  "

$ws.Range("D38").Value = "'48.81"

$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "'  -5.58%  "

$ws.Range("E40").Value = "'  -12.06%  "

$ws.Range("D41").Value = "'8.22"
$ws.Range("E41").Value = "'  -3.43%  "

$ws.Range("D42").Value = "'0.113"
$ws.Range("E42").Value = "'  -3.49%  "

$ws.Range("D43").Value = "'39.30"
$ws.Range("E43").Value = "'  -1.46%  "

$ws.Range("D44").Value = "'0.266"
$ws.Range("E44").Value = "'  -5.47%  "

$ws.Range("D45").Value = "'2.683.40"
$ws.Range("E45").Value = "'  -0.68%  "

$ws.Range("D46").Value = "'133.79"
$ws.Range("E46").Value = "'  -1.05%  "

$ws.Range("E47").Value = "'  -1.54%  "

$ws.Range("D49").Value = "'338.03"
$ws.Range("E49").Value = "'  -7.08%  "

$ws.Range("E50").Value = "'  -1.91%  "

$ws.Range("D51").Value = "'21.39"
$ws.Range("E51").Value = "'  -6.68%  "
